$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'68.741.30"
$ws.Range('E2').Value = "'  +1.39%  "
$ws.Range('D3').Value = "'3.324.31"
$ws.Range('E3').Value = "'  +1.81%  "
$ws.Range('E4').Value = "'  +0.08%  "
$ws.Range('D5').Value = "'188.91"
$ws.Range('E5').Value = "'  +1.99%  "
$ws.Range('D6').Value = "'588.31"
$ws.Range('E6').Value = "'  +1.39%  "
$ws.Range('E7').Value = "'  +0.08%  "
$ws.Range('E8').Value = "'  +0.67%  "
$ws.Range('E9').Value = "'  +1.44%  "
$ws.Range('D10').Value = "'6.73"
$ws.Range('E10').Value = "'  +2.63%  "
$ws.Range('D11').Value = "'0.413"
$ws.Range('E11').Value = "'  +1.13%  "
$ws.Range('D12').Value = "'3.906.95"
$ws.Range('E12').Value = "'  +2.08%  "
$ws.Range('E13').Value = "'  -0.97%  "
$ws.Range('D14').Value = "'27.99"
$ws.Range('E14').Value = "'  +1.98%  "
$ws.Range('D15').Value = "'68.885.95"
$ws.Range('E15').Value = "'  +1.53%  "
$ws.Range('E16').Value = "'  +0.87%  "
$ws.Range('D17').Value = "'3.329.01"
$ws.Range('E17').Value = "'  +1.26%  "
$ws.Range('D18').Value = "'446.58"
$ws.Range('E18').Value = "'  +11.89%  "
$ws.Range('D19').Value = "'5.78"
$ws.Range('D20').Value = "'13.77"
$ws.Range('E20').Value = "'  +2.00%  "
$ws.Range('D21').Value = "'7.77"
$ws.Range('E21').Value = "'  +2.31%  "
$ws.Range('D22').Value = "'75.67"
$ws.Range('E22').Value = "'  +6.58%  "
$ws.Range('E23').Value = "'  -0.11%  "
$ws.Range('B24').Value = "'Polygon"
$ws.Range('C24').Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range('D24').Value = "'0.521"
$ws.Range('E24').Value = "'  +2.15%  "
$ws.Range('B25').Value = "'WrappedeETH"
$ws.Range('C25').Value = "'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range('D25').Value = "'3.496.41"
$ws.Range('E25').Value = "'  +2.51%  "
$ws.Range('E26').Value = "'  +2.47%  "
$ws.Range('E27').Value = "'  +1.78%  "
$ws.Range('D28').Value = "'9.35"
$ws.Range('E28').Value = "'  -1.69%  "
$ws.Range('E29').Value = "'  -0.53%  "
$ws.Range('E30').Value = "'  +2.35%  "
$ws.Range('D31').Value = "'23.14"
$ws.Range('E31').Value = "'  +2.44%  "
$ws.Range('D32').Value = "'5.46"
$ws.Range('E32').Value = "'  -0.26%  "
$ws.Range('E33').Value = "'  +1.66%  "
$ws.Range('D34').Value = "'6.92"
$ws.Range('E34').Value = "'  -0.11%  "
$ws.Range('E35').Value = "'  +0.07%  "
$ws.Range('D36').Value = "'1.54"
$ws.Range('E36').Value = "'  +5.66%  "
$ws.Range('D37').Value = "'163.65"
$ws.Range('E37').Value = "'  +0.85%  "
$ws.Range('E38').Value = "'  +1.99%  "
$ws.Range('D39').Value = "'26.99"
$ws.Range('E39').Value = "'  +0.56%  "
$ws.Range('D40').Value = "'4.56"
$ws.Range('E40').Value = "'  +1.38%  "
$ws.Range('D41').Value = "'0.796"
$ws.Range('E41').Value = "'  -1.30%  "
$ws.Range('E42').Value = "'  +1.90%  "
$ws.Range('D43').Value = "'2.696.31"
$ws.Range('E43').Value = "'  +0.97%  "
$ws.Range('E44').Value = "'  +2.52%  "
$ws.Range('D45').Value = "'41.06"
$ws.Range('E45').Value = "'  +0.98%  "
$ws.Range('D46').Value = "'0.0683"
$ws.Range('E46').Value = "'  +0.42%  "
$ws.Range('D47').Value = "'25.04"
$ws.Range('E47').Value = "'  +1.89%  "
$ws.Range('D48').Value = "'332.02"
$ws.Range('E48').Value = "'  -0.90%  "
$ws.Range('D49').Value = "'0.0281"
$ws.Range('E49').Value = "'  +2.40%  "
$ws.Range('D50').Value = "'32.22"
$ws.Range('E50').Value = "'  +5.39%  "
$ws.Range('E51').Value = "'  +3.48%  "
